$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header label in I1 from "syst_c" to "syst1_c"
$ws.Range("I1").Value = "syst1_c"

# Scale the value/stat_u/syst_u/syst_c columns (F:I) for data rows 2-17.
# Rows 8, 9, 16 and 17 are scaled by 2; all other data rows are scaled by 4.
$doubleRows = @(8, 9, 16, 17)

for ($row = 2; $row -le 17; $row++) {
    if ($doubleRows -contains $row) {
        $factor = 2
    } else {
        $factor = 4
    }
    foreach ($col in @("F", "G", "H", "I")) {
        $addr = $col + $row
        $cur = $ws.Range($addr).Value2
        $ws.Range($addr).Value = $cur * $factor
    }
}

# Move the active selection to G27 (it was I18 before).
$ws.Range("G27").Select()
